$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.387082182999791
$ws.Range("C2").Value = 0.3183208541906311
$ws.Range("D2").Value = 0.02230414481786003
$ws.Range("F2").Value = 1.153346434063565
$ws.Range("G2").Value = 0.002427394107367429
$ws.Range("N2").Value = 1.065748182730892

$ws.Range("B3").Value = 1.246956845320142
$ws.Range("C3").Value = 0.2783575579910007
$ws.Range("D3").Value = 0.02288450579839374
$ws.Range("F3").Value = 1.120154855645808
$ws.Range("G3").Value = 0.002432599193044694
$ws.Range("N3").Value = 1.084394648645083

$ws.Range("B4").Value = 1.161469763897401
$ws.Range("C4").Value = 0.2538598965445544
$ws.Range("D4").Value = 0.0232602923208205
$ws.Range("F4").Value = 1.100782598082063
$ws.Range("G4").Value = 0.002435960986179827
$ws.Range("N4").Value = 1.096411244357649

$ws.Range("B5").Value = 1.126769461237473
$ws.Range("C5").Value = 0.2438862792383816
$ws.Range("D5").Value = 0.02341828002841329
$ws.Range("F5").Value = 1.093139313791085
$ws.Range("G5").Value = 0.002437372794601322
$ws.Range("N5").Value = 1.101450403116342

$ws.Range("B6").Value = 1.121015680594155
$ws.Range("C6").Value = 0.2422307147604954
$ws.Range("D6").Value = 0.02344480571099172
$ws.Range("F6").Value = 1.091885249884768
$ws.Range("G6").Value = 0.002437609756344844
$ws.Range("N6").Value = 1.102295733697219

$ws.Range("B7").Value = 1.161001233959496
$ws.Range("C7").Value = 0.2537253517504325
$ws.Range("D7").Value = 0.02326240340421393
$ws.Range("F7").Value = 1.100678504527338
$ws.Range("G7").Value = 0.002435979856923625
$ws.Range("N7").Value = 1.096478628631438

$ws.Range("B8").Value = 1.338651442927357
$ws.Range("C8").Value = 0.3045327969826133
$ws.Range("D8").Value = 0.02250019546622184
$ws.Range("F8").Value = 1.14169148140742
$ws.Range("G8").Value = 0.002429154486996137
$ws.Range("N8").Value = 1.072059424782573

$ws.Range("B9").Value = 1.69150181053169
$ws.Range("C9").Value = 0.404518578020145
$ws.Range("D9").Value = 0.02116143977373497
$ws.Range("F9").Value = 1.230220691425259
$ws.Range("G9").Value = 0.002417079216009228
$ws.Range("N9").Value = 1.02869546342021

$ws.Range("B10").Value = 1.95365370078099
$ws.Range("C10").Value = 0.4782496078823328
$ws.Range("D10").Value = 0.02027519159223612
$ws.Range("F10").Value = 1.300362191804169
$ws.Range("G10").Value = 0.002408996290415925
$ws.Range("N10").Value = 0.9996171040946376

$ws.Range("B11").Value = 2.073587468012079
$ws.Range("C11").Value = 0.5118635238475235
$ws.Range("D11").Value = 0.01989364969560281
$ws.Range("F11").Value = 1.333413815362775
$ws.Range("G11").Value = 0.002405488426891303
$ws.Range("N11").Value = 0.9869978592314776

$ws.Range("B12").Value = 2.119103845253392
$ws.Range("C12").Value = 0.524603681369797
$ws.Range("D12").Value = 0.01975232162983431
$ws.Range("F12").Value = 1.346096830252151
$ws.Range("G12").Value = 0.00240418425299315
$ws.Range("N12").Value = 0.9823072994363713

$ws.Range("B13").Value = 2.109296590859628
$ws.Range("C13").Value = 0.5218593458553187
$ws.Range("D13").Value = 0.01978261817443183
$ws.Range("F13").Value = 1.343357843459785
$ws.Range("G13").Value = 0.002404464057161325
$ws.Range("N13").Value = 0.9833135693826538

$ws.Range("B14").Value = 2.077330108756144
$ws.Range("C14").Value = 0.5129114334778251
$ws.Range("D14").Value = 0.01988195908731782
$ws.Range("F14").Value = 1.334453891106449
$ws.Range("G14").Value = 0.002405380648057408
$ws.Range("N14").Value = 0.9866101957170343

$ws.Range("B15").Value = 2.057762816958814
$ws.Range("C15").Value = 0.5074320726979522
$ws.Range("D15").Value = 0.01994322030413898
$ws.Range("F15").Value = 1.329021799644579
$ws.Range("G15").Value = 0.00240594523075575
$ws.Range("N15").Value = 0.988640960963755

$ws.Range("B16").Value = 1.945829620588199
$ws.Range("C16").Value = 0.4760543987545702
$ws.Range("D16").Value = 0.0203005650751491
$ws.Range("F16").Value = 1.298225429489577
$ws.Range("G16").Value = 0.002409228928693218
$ws.Range("N16").Value = 1.000454095877165

$ws.Range("B17").Value = 1.877337932089574
$ws.Range("C17").Value = 0.456824519846009
$ws.Range("D17").Value = 0.02052535294505109
$ws.Range("F17").Value = 1.279627668674294
$ws.Range("G17").Value = 0.002411286582599758
$ws.Range("N17").Value = 1.007857349808079

$ws.Range("B18").Value = 1.838007144351991
$ws.Range("C18").Value = 0.4457708553816246
$ws.Range("D18").Value = 0.02065667691800588
$ws.Range("F18").Value = 1.26903824145954
$ws.Range("G18").Value = 0.002412486015664983
$ws.Range("N18").Value = 1.012172760759583

$ws.Range("B19").Value = 1.82470130600467
$ws.Range("C19").Value = 0.4420294287715478
$ws.Range("D19").Value = 0.02070148858447407
$ws.Range("F19").Value = 1.265471230347643
$ws.Range("G19").Value = 0.002412894862193427
$ws.Range("N19").Value = 1.01364370172212

$ws.Range("B20").Value = 1.884622367321754
$ws.Range("C20").Value = 0.4588708585620793
$ws.Range("D20").Value = 0.02050121327106602
$ws.Range("F20").Value = 1.281596286217081
$ws.Range("G20").Value = 0.00241106589492569
$ws.Range("N20").Value = 1.007063331142108

$ws.Range("B21").Value = 2.086716704835055
$ws.Range("C21").Value = 0.5155393395431815
$ws.Range("D21").Value = 0.01985269430243974
$ws.Range("F21").Value = 1.337064642384888
$ws.Range("G21").Value = 0.0024051107678359
$ws.Range("N21").Value = 0.9856395018891142

$ws.Range("B22").Value = 2.21938140196437
$ws.Range("C22").Value = 0.5526416773549272
$ws.Range("D22").Value = 0.01944725257777602
$ws.Range("F22").Value = 1.37429152735001
$ws.Range("G22").Value = 0.00240135961168756
$ws.Range("N22").Value = 0.9721514191360203

$ws.Range("B23").Value = 2.148521499404637
$ws.Range("C23").Value = 0.532833149412852
$ws.Range("D23").Value = 0.01966194548316658
$ws.Range("F23").Value = 1.354332751996878
$ws.Range("G23").Value = 0.002403348830829871
$ws.Range("N23").Value = 0.9793030807515279

$ws.Range("B24").Value = 1.88132893106075
$ws.Range("C24").Value = 0.4579457031262564
$ws.Range("D24").Value = 0.02051212030954908
$ws.Range("F24").Value = 1.280705954617375
$ws.Range("G24").Value = 0.002411165616465798
$ws.Range("N24").Value = 1.007422122685034

$ws.Range("B25").Value = 1.595548276855197
$ws.Range("C25").Value = 0.3774262368228847
$ws.Range("D25").Value = 0.0215066826389041
$ws.Range("F25").Value = 1.205387353482024
$ws.Range("G25").Value = 0.002420206703768194
$ws.Range("N25").Value = 1.03993996558367

